$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with new columns P and Q, copying formatting from O1
$hdr = $ws.Range("P1:Q1")
$hdr.Borders.LineStyle = 1       # thin border, all sides
$hdr.Borders.Weight = 2
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108 # xlCenter
$hdr.VerticalAlignment = -4160   # xlTop

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

for ($r = 2; $r -le 25; $r++) {
    # Flip values in existing columns
    $ws.Cells.Item($r, 9).Value = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: 2 -> 1

    # New columns P and Q
    $ws.Cells.Item($r, 16).Value = 2  # P
    $ws.Cells.Item($r, 17).Value = 2  # Q
}

$wb.Save()
